$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.038.30'
$ws.Range("E2").Value = '  +5.26%  '
$ws.Range("D3").Value = '2.417.57'
$ws.Range("E3").Value = '  +1.86%  '
$ws.Range("E4").Value = '  +0.46%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.05'
$ws.Range("E5").Value = '  +2.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.23'
$ws.Range("E6").Value = '  +6.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.39%  '
$ws.Range("E8").Value = '  +2.28%  '
$ws.Range("D9").Value = '2.450.83'
$ws.Range("E9").Value = '  +3.38%  '
$ws.Range("E10").Value = '  +5.83%  '
$ws.Range("E11").Value = '  +0.71%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.25'
$ws.Range("E12").Value = '  +3.06%  '
$ws.Range("E13").Value = '  +4.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.36'
$ws.Range("E14").Value = '  +6.59%  '
$ws.Range("E15").Value = '  +8.10%  '
$ws.Range("D16").Value = '2.858.16'
$ws.Range("D17").Value = '62.775.21'
$ws.Range("E17").Value = '  +4.93%  '
$ws.Range("D18").Value = '2.445.46'
$ws.Range("E18").Value = '  +2.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.93'
$ws.Range("E19").Value = '  -1.90%  '
$ws.Range("E20").Value = '  +4.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '329.06'
$ws.Range("E21").Value = '  +2.33%  '
$ws.Range("E22").Value = '  +2.19%  '
$ws.Range("E23").Value = '  +14.21%  '
$ws.Range("E24").Value = '  -0.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.62'
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '635.78'
$ws.Range("E26").Value = '  +14.06%  '
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.53'
$ws.Range("E27").Value = '  +4.86%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0993'
$ws.Range("E28").Value = '  +7.58%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("C29").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D29").Value = '2.538.18'
$ws.Range("E29").Value = '  +2.01%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.24'
$ws.Range("E30").Value = '  +2.44%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.42'
$ws.Range("E31").Value = '  +9.22%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.139'
$ws.Range("E32").Value = '  +6.15%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.85'
$ws.Range("E33").Value = '  +3.84%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.50'
$ws.Range("E34").Value = '  +4.55%  '
$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.995'
$ws.Range("E35").Value = '  -0.46%  '
$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.76'
$ws.Range("E36").Value = '  +4.90%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.375'
$ws.Range("E37").Value = '  +2.45%  '
$ws.Range("B38").Value = 'Monero'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '153.02'
$ws.Range("E38").Value = '  -0.10%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("E39").Value = '  +8.82%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.80'
$ws.Range("E40").Value = '  +16.03%  '
$ws.Range("E41").Value = '  +3.03%  '
$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.77'
$ws.Range("E42").Value = '  +7.92%  '
$ws.Range("B43").Value = 'USDe'
$ws.Range("C43").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("B44").Value = 'BabyDogeCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D44").Value = '0.0₆0286'
$ws.Range("E44").Value = '  -4.10%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '145.26'
$ws.Range("E45").Value = '  +3.92%  '
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.61'
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.53'
$ws.Range("E47").Value = '  +7.84%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.603'
$ws.Range("E48").Value = '  +3.08%  '
$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0518'
$ws.Range("E49").Value = '  +3.61%  '
$ws.Range("B50").Value = 'WhiteBITCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '13.05'
$ws.Range("E50").Value = '  +11.68%  '
$ws.Range("B51").Value = 'Stellar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0920'
$ws.Range("E51").Value = '  +2.37%  '
